# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gains a new column (between the existing
# "In Advance" column and the "Late" column) so a Variable-Instalment loan
# repayment schedule can report an extra figure. Inserting the column
# shifts the old "Late" / "Outstanding" columns one slot to the right,
# which is exactly what Excel's own Insert-a-column command does, and
# also grows the sheet's column-width/dimension bookkeeping accordingly.
# The sheet also becomes the active tab/selection, matching how the
# workbook was left after the edit.

$wb = $excel.ActiveWorkbook

$wsSched  = $wb.Worksheets.Item("Repayment schedule")

# --- NewLoanInput sheet: no longer the tab that's selected/active ---
# (handled implicitly below once we activate "Repayment schedule" and
#  select a cell on it -- Excel moves tabSelected/activeTab bookkeeping
#  to whichever sheet/range is active last.)

# --- Repayment schedule sheet: insert new column N ---
# Existing layout (1-indexed columns):
#   ... L=Late-Due M=Late-Paid N=In Advance O=Late P=Outstanding
# After inserting a blank column at 14 (N), the old N/O/P shift right to
# O/P/Q, leaving the new N blank and ready for its own heading/values.
$wsSched.Columns.Item(14).Insert()

# Match the width Excel applied to the freshly inserted column (same
# width as its left neighbour, column M).
$wsSched.Columns.Item(14).ColumnWidth = $wsSched.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with the same cell
# selected as in the edited workbook.
$wsSched.Activate()
$wsSched.Range("I18").Select()
